$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.489.82"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "2.986.38"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.65"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.00"
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.592"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.68"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "3.455.01"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.46"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.79"
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").Value = "2.987.92"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.23"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "51.483.19"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.11"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.58"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.27"
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.43"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("E25").Value = "  +3.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.85"
$ws.Range("E26").Value = "  -3.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.36"
$ws.Range("E27").Value = "  -2.70%  "
$ws.Range("E28").Value = "  +2.69%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.09"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("E32").Value = "  +3.56%  "
$ws.Range("E33").Value = "  +3.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.41"
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0444"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.29"
$ws.Range("E38").Value = "  +4.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.93"
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("E40").Value = "  +3.97%  "
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.85"
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.81"
$ws.Range("E43").Value = "  +12.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "123.24"
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.40"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("D49").Value = "2.027.21"
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("D50").Value = "3.281.66"
$ws.Range("E50").Value = "  +1.38%  "
$ws.Range("E51").Value = "  +1.41%  "
